$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (was date 44400 / old-row-9 style data)
$ws.Range("D2").Value = 44397
$ws.Range("M2").Value = 60
$ws.Range("N2").Value = 11000
$ws.Range("O2").Value = 11000
$ws.Range("P2").Value = 11000
$ws.Range("S2").Value = 786

# Row 3 (swap with row 4)
$ws.Range("D3").Value = 44176
$ws.Range("M3").Value = 250

# Row 4 (swap with row 3)
$ws.Range("D4").Value = 44162
$ws.Range("M4").Value = 120

# Row 5 (swap with row 6)
$ws.Range("D5").Value = 44309
$ws.Range("M5").Value = 300
$ws.Range("N5").Value = 7000
$ws.Range("O5").Value = 7000
$ws.Range("P5").Value = 7000
$ws.Range("S5").Value = 500

# Row 6 (swap with row 5)
$ws.Range("D6").Value = 44208
$ws.Range("M6").Value = 210
$ws.Range("N6").Value = 10000
$ws.Range("O6").Value = 10000
$ws.Range("P6").Value = 10000
$ws.Range("S6").Value = 714

# Row 8 (was date 44491 / old-row-2 style data)
$ws.Range("D8").Value = 44400
$ws.Range("M8").Value = 100
$ws.Range("N8").Value = 10000
$ws.Range("O8").Value = 10000
$ws.Range("P8").Value = 10000
$ws.Range("Q8").Value = "$/caja 14 kilos"
$ws.Range("S8").Value = 714

# Row 9 (was date 44397 / old-row-8 style data)
$ws.Range("D9").Value = 44491
$ws.Range("M9").Value = 180
$ws.Range("N9").Value = 9000
$ws.Range("O9").Value = 9000
$ws.Range("P9").Value = 9000
$ws.Range("Q9").Value = "$/caja 14 kilos empedrada"
$ws.Range("S9").Value = 643
